# The deck's single slide master currently carries the "Integral" theme
# (ppt/theme/theme2.xml) while the notes master carries the default
# "Office Theme" (ppt/theme/theme1.xml). The authored change swaps the
# two themes' color content so the slide master's theme becomes the
# "Office Theme" palette. We reproduce that by writing the 12 DrawingML
# theme colors (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) of the
# presentation's active theme to the stock Office Theme RGB values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# RGB() packs as R + G*256 + B*65536, matching the VBA RGB() convention.
function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Index -> element : 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$tcs.Colors(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1      000000
$tcs.Colors(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Colors(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2      44546A
$tcs.Colors(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Colors(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Colors(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Colors(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Colors(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Colors(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Colors(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Colors(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Colors(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink 954F72
